# Auto-generated edit script applying market-data refresh to Sargatanas_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across all 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 16666.334
$ws.Range("I21").Value = 16666.334
$ws.Range("K21").Value = 16666.334
$ws.Range("M21").Value = -16198.334
$ws.Range("H23").Value = 16666.334
$ws.Range("I23").Value = 16666.334
$ws.Range("K23").Value = 16666.334
$ws.Range("M23").Value = -16432.334
$ws.Range("H107").Value = 5000367.5
$ws.Range("I107").Value = 5682108.5
$ws.Range("K107").Value = 5682108.5
$ws.Range("M107").Value = -5680188.5
$ws.Range("H115").Value = 439.375
$ws.Range("I115").Value = 439.375
$ws.Range("K115").Value = 1318.125
$ws.Range("M115").Value = 248.875
$ws.Range("H132").Value = 2435.75
$ws.Range("I132").Value = 2415.9678
$ws.Range("J132").Value = 2640.1667
$ws.Range("K132").Value = 7247.903399999999
$ws.Range("L132").Value = 7920.500100000001
$ws.Range("M132").Value = -4717.903399999999
$ws.Range("N132").Value = -12980.5001
$ws.Range("H137").Value = 5654.0356
$ws.Range("I137").Value = 4490.1055
$ws.Range("J137").Value = 8111.222
$ws.Range("K137").Value = 13470.3165
$ws.Range("L137").Value = 24333.666
$ws.Range("M137").Value = -10920.3165
$ws.Range("N137").Value = -29433.666
$ws.Range("H138").Value = 3363.4768
$ws.Range("I138").Value = 1260.2094
$ws.Range("J138").Value = 5466.744
$ws.Range("K138").Value = 3780.6282
$ws.Range("L138").Value = 16400.232
$ws.Range("M138").Value = 1359.3718
$ws.Range("N138").Value = -26680.232
$ws.Range("H141").Value = 1929.25
$ws.Range("I141").Value = 187.11111
$ws.Range("K141").Value = 561.3333299999999
$ws.Range("M141").Value = 4618.666670000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3993.8052
$ws.Range("I32").Value = 3388.2764
$ws.Range("K32").Value = 3388.2764
$ws.Range("M32").Value = -3101.2764
$ws.Range("H61").Value = 33342276
$ws.Range("I61").Value = 5657
$ws.Range("J61").Value = 71441270
$ws.Range("K61").Value = 5657
$ws.Range("L61").Value = 71441270
$ws.Range("M61").Value = -5445
$ws.Range("N61").Value = -71441694
$ws.Range("H88").Value = 1550
$ws.Range("I88").Value = 900
$ws.Range("J88").Value = 1642.8572
$ws.Range("K88").Value = 900
$ws.Range("L88").Value = 1642.8572
$ws.Range("M88").Value = -494
$ws.Range("N88").Value = -2454.8572
$ws.Range("H91").Value = 1550
$ws.Range("I91").Value = 900
$ws.Range("J91").Value = 1642.8572
$ws.Range("K91").Value = 900
$ws.Range("L91").Value = 1642.8572
$ws.Range("M91").Value = 504
$ws.Range("N91").Value = -4450.8572
$ws.Range("H110").Value = 83333840
$ws.Range("I110").Value = 667
$ws.Range("K110").Value = 667
$ws.Range("M110").Value = 1378
$ws.Range("H136").Value = 33342276
$ws.Range("I136").Value = 5657
$ws.Range("J136").Value = 71441270
$ws.Range("K136").Value = 16971
$ws.Range("L136").Value = 214323810
$ws.Range("M136").Value = -14421
$ws.Range("N136").Value = -214328910

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 7333
$ws.Range("I31").Value = 6999.5
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 6999.5
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -6747.5
$ws.Range("N31").Value = -8504
$ws.Range("H105").Value = 3170.3333
$ws.Range("I105").Value = 2184.0715
$ws.Range("J105").Value = 5142.857
$ws.Range("K105").Value = 2184.0715
$ws.Range("L105").Value = 5142.857
$ws.Range("M105").Value = -437.0715
$ws.Range("N105").Value = -8636.857
$ws.Range("H109").Value = 53020
$ws.Range("J109").Value = 53020
$ws.Range("L109").Value = 53020
$ws.Range("N109").Value = -55794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4679.6875
$ws.Range("I16").Value = 837.8570999999999
$ws.Range("J16").Value = 7667.778
$ws.Range("K16").Value = 837.8570999999999
$ws.Range("L16").Value = 7667.778
$ws.Range("M16").Value = -550.8570999999999
$ws.Range("N16").Value = -8241.778
$ws.Range("H31").Value = 8555833
$ws.Range("I31").Value = 2333.1667
$ws.Range("J31").Value = 10111015
$ws.Range("K31").Value = 2333.1667
$ws.Range("L31").Value = 10111015
$ws.Range("M31").Value = -2038.1667
$ws.Range("N31").Value = -10111605
$ws.Range("H34").Value = 8555833
$ws.Range("I34").Value = 2333.1667
$ws.Range("J34").Value = 10111015
$ws.Range("K34").Value = 2333.1667
$ws.Range("L34").Value = 10111015
$ws.Range("M34").Value = -2131.1667
$ws.Range("N34").Value = -10111419
$ws.Range("H39").Value = 7590.4
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 7590.4
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 7590.4
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -8372.4
$ws.Range("H41").Value = 47500
$ws.Range("I41").Value = 30000
$ws.Range("K41").Value = 30000
$ws.Range("M41").Value = -29572
$ws.Range("H42").Value = 26625
$ws.Range("I42").Value = 27000
$ws.Range("J42").Value = 26250
$ws.Range("K42").Value = 27000
$ws.Range("L42").Value = 26250
$ws.Range("M42").Value = -26407
$ws.Range("N42").Value = -27436
$ws.Range("H49").Value = 7590.4
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 7590.4
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 7590.4
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = -7954.4
$ws.Range("H58").Value = 6925.816
$ws.Range("I58").Value = 2360.5
$ws.Range("J58").Value = 11034.6
$ws.Range("K58").Value = 2360.5
$ws.Range("L58").Value = 11034.6
$ws.Range("M58").Value = -2157.5
$ws.Range("N58").Value = -11440.6
$ws.Range("H75").Value = 43000
$ws.Range("J75").Value = 43000
$ws.Range("L75").Value = 43000
$ws.Range("N75").Value = -44996
$ws.Range("H78").Value = 43000
$ws.Range("J78").Value = 43000
$ws.Range("L78").Value = 129000
$ws.Range("N78").Value = -138984
$ws.Range("H107").Value = 1314.825
$ws.Range("I107").Value = 767.05
$ws.Range("J107").Value = 1862.6
$ws.Range("K107").Value = 767.05
$ws.Range("L107").Value = 1862.6
$ws.Range("M107").Value = 1152.95
$ws.Range("N107").Value = -5702.6
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H109").Value = 59376
$ws.Range("J109").Value = 59376
$ws.Range("L109").Value = 59376
$ws.Range("N109").Value = -61456
$ws.Range("H113").Value = 4679.6875
$ws.Range("I113").Value = 837.8570999999999
$ws.Range("J113").Value = 7667.778
$ws.Range("K113").Value = 837.8570999999999
$ws.Range("L113").Value = 7667.778
$ws.Range("M113").Value = 1332.1429
$ws.Range("N113").Value = -12007.778
$ws.Range("H134").Value = 5820.4565
$ws.Range("I134").Value = 2269.8
$ws.Range("J134").Value = 10047.429
$ws.Range("K134").Value = 6809.400000000001
$ws.Range("L134").Value = 30142.287
$ws.Range("M134").Value = -4274.400000000001
$ws.Range("N134").Value = -35212.287
$ws.Range("H136").Value = 6925.816
$ws.Range("I136").Value = 2360.5
$ws.Range("J136").Value = 11034.6
$ws.Range("K136").Value = 7081.5
$ws.Range("L136").Value = 33103.8
$ws.Range("M136").Value = -4531.5
$ws.Range("N136").Value = -38203.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 3077.0833
$ws.Range("I98").Value = 143.33333
$ws.Range("J98").Value = 4055
$ws.Range("K98").Value = 429.99999
$ws.Range("L98").Value = 12165
$ws.Range("M98").Value = 1068.00001
$ws.Range("N98").Value = -15161
$ws.Range("H132").Value = 7029.8
$ws.Range("I132").Value = 2410.8572
$ws.Range("K132").Value = 21697.7148
$ws.Range("M132").Value = -19167.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 800652
$ws.Range("I107").Value = 1000790.1
$ws.Range("J107").Value = 99.5
$ws.Range("K107").Value = 1000790.1
$ws.Range("L107").Value = 99.5
$ws.Range("M107").Value = -998870.1
$ws.Range("N107").Value = -3939.5
$ws.Range("H126").Value = 12400
$ws.Range("J126").Value = 12400
$ws.Range("L126").Value = 37200
$ws.Range("N126").Value = -42140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6177256.5
$ws.Range("I132").Value = 10640388
$ws.Range("J132").Value = 7634.147
$ws.Range("K132").Value = 31921164
$ws.Range("L132").Value = 22902.441
$ws.Range("M132").Value = -31918634
$ws.Range("N132").Value = -27962.441
$ws.Range("H136").Value = 6090.2354
$ws.Range("I136").Value = 2113.0977
$ws.Range("J136").Value = 12129.593
$ws.Range("K136").Value = 6339.293099999999
$ws.Range("L136").Value = 36388.779
$ws.Range("M136").Value = -3789.293099999999
$ws.Range("N136").Value = -41488.779

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 10000
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10460
$ws.Range("H109").Value = 59376
$ws.Range("J109").Value = 59376
$ws.Range("L109").Value = 59376
$ws.Range("N109").Value = -62150
$ws.Range("H122").Value = 195374.72
$ws.Range("I122").Value = 289212.06
$ws.Range("J122").Value = 7700
$ws.Range("K122").Value = 867636.1799999999
$ws.Range("L122").Value = 23100
$ws.Range("M122").Value = -865186.1799999999
$ws.Range("N122").Value = -28000
$ws.Range("H132").Value = 12823397
$ws.Range("I132").Value = 15154126
$ws.Range("J132").Value = 4389.3335
$ws.Range("K132").Value = 45462378
$ws.Range("L132").Value = 13168.0005
$ws.Range("M132").Value = -45459848
$ws.Range("N132").Value = -18228.0005
